# Applies the "Issues fixes and reports" edit:
#  - AMSIN sheet: append rows 33 & 34 (new registration-history records)
#  - AMS sheet:   fix formatting on row 30 (+tiny B30 value correction),
#                 append rows 31 & 32 (new registration-history records)

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $text) {
    # Force text storage (avoids Excel's automatic date/number parsing of
    # strings like "2023-02-17"), then drop back to a plain/General style
    # so the cell isn't left tagged with a lingering "@" number format.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

function Set-NumberCell($ws, $addr, $num) {
    $ws.Range($addr).Value = $num
    $ws.Range($addr).Style = "Normal"
}

function Set-DateTimeCell($ws, $addr, $num) {
    $ws.Range($addr).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range($addr).Value = $num
}

# ---------------------------------------------------------------------
# Sheet "AMSIN": add rows 33 and 34
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Set-TextCell     $wsAmsin "A33" "2023-02-17"
Set-DateTimeCell $wsAmsin "B33" 44974.46063145834
Set-TextCell     $wsAmsin "C33" "pay173fstccycle"
Set-NumberCell   $wsAmsin "D33" 44
Set-NumberCell   $wsAmsin "E33" 44
Set-NumberCell   $wsAmsin "F33" 0
Set-NumberCell   $wsAmsin "G33" 1.25

Set-TextCell     $wsAmsin "A34" "2023-02-20"
Set-DateTimeCell $wsAmsin "B34" 44977.43036248843
Set-TextCell     $wsAmsin "C34" "173payflow"
Set-NumberCell   $wsAmsin "D34" 44
Set-NumberCell   $wsAmsin "E34" 44
Set-NumberCell   $wsAmsin "F34" 0
Set-NumberCell   $wsAmsin "G34" 1.44

# ---------------------------------------------------------------------
# Sheet "AMS": normalize styling on row 30, correct B30, add rows 31/32
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Range("A30").Style = "Normal"
$wsAms.Range("C30").Style = "Normal"
$wsAms.Range("D30").Style = "Normal"
$wsAms.Range("E30").Style = "Normal"
$wsAms.Range("F30").Style = "Normal"
$wsAms.Range("G30").Style = "Normal"
$wsAms.Range("B30").Value = 44964.5773160301

Set-TextCell     $wsAms "A31" "2023-02-20"
Set-DateTimeCell $wsAms "B31" 44977.60840824074
Set-TextCell     $wsAms "C31" "173payflow"
Set-NumberCell   $wsAms "D31" 44
Set-NumberCell   $wsAms "E31" 44
Set-NumberCell   $wsAms "F31" 0
Set-NumberCell   $wsAms "G31" 1.49

# Row 32 matches the file's "un-styled" look (no explicit style, like the
# pre-fix row 30) -- plain value writes already default to that; the two
# text cells need the date-parse dodge but then get their format stripped
# back to the implicit/no-style state with ClearFormats.
$wsAms.Range("A32").NumberFormat = "@"
$wsAms.Range("A32").Value = "2023-02-20"
$wsAms.Range("A32").ClearFormats()
$wsAms.Range("B32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Range("B32").Value = 44977.82336558346
$wsAms.Range("C32").NumberFormat = "@"
$wsAms.Range("C32").Value = "173livepay"
$wsAms.Range("C32").ClearFormats()
$wsAms.Range("D32").Value = 44
$wsAms.Range("E32").Value = 44
$wsAms.Range("F32").Value = 0
$wsAms.Range("G32").Value = 1.35
